$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New snippet rows describing the sessionData APIs (Outlook), appended
# to the bottom of the "Snippets" table.
$rows = @(
    @("AppointmentCompose", "sessionData", $null, "outlook-session-data-apis", "getAllSessionData"),
    @("MessageCompose",     "sessionData", $null, "outlook-session-data-apis", "getAllSessionData"),
    @("SessionData",        "setAsync",    1,     "outlook-session-data-apis", "setSessionData"),
    @("SessionData",        "getAsync",    1,     "outlook-session-data-apis", "getSessionData"),
    @("SessionData",        "getAllAsync", 1,     "outlook-session-data-apis", "getAllSessionData"),
    @("SessionData",        "removeAsync", 1,     "outlook-session-data-apis", "removeSessionData"),
    @("SessionData",        "clearAsync",  1,     "outlook-session-data-apis", "clearSessionData")
)

$startRow = $lo.Range.Row + $lo.Range.Rows.Count
$r = $startRow

foreach ($row in $rows) {
    # Insert a fresh row at the bottom of the table, copying formatting
    # (style) down from the row immediately above it.
    $ws.Rows.Item($r).Insert(-4121, -4163)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $r = $r + 1
}

$lastRow = $r - 1
$null = $lo.Resize($ws.Range($ws.Cells.Item($lo.Range.Row, 1), $ws.Cells.Item($lastRow, 5)))

$null = $ws.Range("E" + $lastRow).Select()
